$wb = $excel.ActiveWorkbook

# zh-cn sheet: update the "f8b7a46d..." row's Correspond Handoff / Handback datetimes
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-22 12:03:54"
$wsZhCn.Range("H4").Value = "2016-03-22 12:04:50"

# de-de sheet: update the "f8b7a46d..." row's Correspond Handoff / Handback datetimes
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-22 12:04:02"
$wsDeDe.Range("H4").Value = "2016-03-22 12:05:04"
